$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed crypto "Price" (D) / "Volume(1h)" (E) figures from the latest
# scrape. Values are numeric-looking text (e.g. "256.62", "0.09%"), so we
# prefix them with a quote to force Excel to store them as literal text
# -- matching the source sheet, which keeps these columns as plain strings
# rather than real numbers/percentages.
$updates = @{
    'D2' = '256.62'
    'E2' = '0.09%'
    'E3' = '-0.01%'
    'D4' = '4.764'
    'E4' = '0.12%'
    'E5' = '-0.27%'
    'D6' = '6.611'
    'E6' = '-0.71%'
    'D7' = '0.8512'
    'E7' = '-1.85%'
    'D8' = '0.9252'
    'E8' = '-1.49%'
    'E9' = '-1.55%'
    'D10' = '0.04169'
    'E10' = '11.11%'
    'D11' = '0.07023'
    'E11' = '-1.14%'
    'E12' = '-3.55%'
    'D13' = '0.09104'
    'E13' = '-1.53%'
    'E14' = '-1.29%'
    'D15' = '0.0006028'
    'E15' = '-94.23%'
    'D16' = '0.006014'
    'E16' = '0.03%'
    'D17' = '3.471'
    'E17' = '-0.46%'
    'D18' = '3.167'
    'E18' = '-1.12%'
    'D19' = '2.199'
    'E19' = '-1.85%'
    'D20' = '0.3042'
    'E20' = '-3.55%'
    'E21' = '0.33%'
    'D22' = '3.897'
    'E22' = '2.35%'
    'D23' = '0.04232'
    'E23' = '0.71%'
    'D24' = '0.001218'
    'E24' = '-0.16%'
    'E25' = '-0.22%'
    'D26' = '0.0001199'
    'E26' = '-0.02%'
    'E27' = '1.94%'
    'D40' = '0.03782'
    'E40' = '-1.03%'
    'D41' = '0.006265'
    'E41' = '1.60%'
    'E42' = '-0.33%'
    'D43' = '0.002198'
    'E43' = '-0.02%'
    'D44' = '0.01410'
    'E44' = '27.11%'
    'D45' = '0.00005345'
    'E45' = '-2.77%'
    'D46' = '0.00000000749'
    'E46' = '-0.02%'
    'D47' = '0.04399'
    'E47' = '-50.28%'
    'E48' = '10,025.11%'
    'D49' = '0.00002098'
    'E49' = '-0.02%'
    'D50' = '0.0001998'
    'E50' = '-0.02%'
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = "'" + $updates[$cellRef]
}
